$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target final state of the schedule table (A1:F17).
# Rows 2-8 keep their original time labels; a new 12:20 lunch slot is
# inserted, pushing every later time slot down by one row and adding two
# new empty slots (17:30, 18:20) at the bottom.

$data = @{
    3  = @("7:50",  "-", "-", "-", "ELT-2A-Circuitos Elétricos 2", "-")
    4  = @("8:40",  "-", "-", "MCT-2A-Circuitos Elétricos 2", "-", "-")
    6  = @("9:50",  "-", "-", "MCT-2A-Circuitos Elétricos 2", "MCT-2A-Programação", "-")
    7  = @("10:40", "ELT-2A-Circuitos Elétricos 2", "ELT-2A-Circuitos Elétricos 2", "-", "-", "-")
    8  = @("11:30", "-", "-", "-", "-", "-")
    9  = @("12:20", "Almoço", "Almoço", "Almoço", "Almoço", "Almoço")
    10 = @("13:00", "-", "-", "-", "-", "-")
    11 = @("13:50", "-", "-", "-", "-", "-")
    12 = @("14:40", "-", "-", "-", "-", "-")
    13 = @("15:30", "Intervalo", "Intervalo", "Intervalo", "Intervalo", "Intervalo")
    14 = @("15:50", "-", "-", "-", "-", "-")
    15 = @("16:40", "-", "-", "-", "-", "-")
    16 = @("17:30", "-", "-", "-", "-", "-")
    17 = @("18:20", "", "", "", "", "")
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
